# Apply the "output-identifier-type" update:
#  - add a new "Version" column at A (shifting Code/Description/Definition one
#    column to the right, B/C/D)
#  - add a new empty "Guide.for.Use" column at E
#  - populate the new Version column with the literal text "1.0" in every
#    data row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell used to force a literal value (e.g. "1.0") to be written as
# TEXT instead of being auto-coerced to a number by the normal Value setter.
# We build it as a formula returning a string, copy it, and paste-special
# "values only" into the real target -- this keeps the destination cell a
# plain shared-string cell (t="s") with no extra number-format / style.
$staging = $ws.Range("ZZ1000")
function Set-TextValue($range, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $staging.Formula = '="' + $escaped + '"'
    $staging.Copy()
    $range.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
    $excel.CutCopyMode = $false
}

# Materialise a present-but-empty cell (<c r="X"/> with no value and no
# style) by touching a formatting property that already matches this
# workbook's default (no fill pattern) -- this forces the cell to exist
# without allocating a new style index.
function Set-EmptyCell($range) {
    $range.Interior.Pattern = [Microsoft.Office.Interop.Excel.XlPattern]::xlPatternNone
}

# ---- Row 1 : header ----
$ws.Range("A1").Value = "Version"
$ws.Range("B1").Value = "Code"
$ws.Range("C1").Value = "Description"
$ws.Range("D1").Value = "Definition"
$ws.Range("E1").Value = "Guide.for.Use"

# ---- Data rows 2-9 ----
$data = @(
    @{ Row = 2; Code = 100; Desc = "O_ISBN10"; Def = "A unique numeric commercial book identifier, issued prior to 2007" },
    @{ Row = 3; Code = 200; Desc = "O_ISBN13"; Def = "A unique numeric commercial book identifier, issued from 2007 onwards" },
    @{ Row = 4; Code = 300; Desc = "O_DOI"; Def = "A unique alphanumeric string assigned to identify content and provide a persistent link to its location on the internet" },
    @{ Row = 5; Code = 400; Desc = "O_URL"; Def = "An address to the location of the output on the internet" },
    @{ Row = 6; Code = 500; Desc = "O_ISSN"; Def = "A serial number used to uniquely identify a serial print publication" },
    @{ Row = 7; Code = 600; Desc = "O_EISSN"; Def = "A serial number used to uniquely identify a serial electronic publication" },
    @{ Row = 8; Code = 700; Desc = "O_PII"; Def = "A unique identifier used by some scientific journals to identify documents" },
    @{ Row = 9; Code = 999; Desc = "Output identifier not available"; Def = "Output identifier not available" }
)

foreach ($item in $data) {
    $r = $item.Row
    Set-TextValue $ws.Range("A$r") "1.0"
    $ws.Range("B$r").Value = $item.Code
    $ws.Range("C$r").Value = $item.Desc
    $ws.Range("D$r").Value = $item.Def
    Set-EmptyCell $ws.Range("E$r")
}

# Clean up the scratch cell used for text-forcing.
$staging.ClearContents()
